$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply AutoFilter on column C ("Grouping") for "uSD" while rows still have
# their original grouping so that the "un-needed" rows remain visible.
$ws.Range("A1:C113").AutoFilter(3, @("uSD"), 7)

# Clear the "un-needed" SDIO pin rows (B, C, D columns) - rows 26, 27, 40, 41
$ws.Range("B26:D26").ClearContents()
$ws.Range("B27:D27").ClearContents()
$ws.Range("B40:D40").ClearContents()
$ws.Range("B41:D41").ClearContents()

# Update selection
$ws.Range("A26").Select()
